$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "86.600.75"
$ws.Range("E2").Value = "  +3.04%  "

$ws.Range("D3").Value = "3.262.74"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'210.89"
$ws.Range("E5").Value = "  -3.57%  "

$ws.Range("D6").Value = "'625.03"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "'0.369"
$ws.Range("E7").Value = "  +19.10%  "

$ws.Range("D8").Value = "'0.682"
$ws.Range("E8").Value = "  +15.66%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").Value = "3.260.08"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("D11").Value = "'0.570"
$ws.Range("E11").Value = "  -5.28%  "

$ws.Range("E12").Value = "  +8.15%  "

$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  -8.62%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.870.26"
$ws.Range("E14").Value = "  +1.28%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'33.95"
$ws.Range("E15").Value = "  +4.40%  "

$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  -1.43%  "

$ws.Range("D17").Value = "86.465.43"
$ws.Range("E17").Value = "  +3.56%  "

$ws.Range("D18").Value = "3.269.47"
$ws.Range("E18").Value = "  +0.95%  "

$ws.Range("D19").Value = "'13.95"
$ws.Range("E19").Value = "  -2.64%  "

$ws.Range("E20").Value = "  -5.93%  "

$ws.Range("D21").Value = "'430.02"
$ws.Range("E21").Value = "  -4.04%  "

$ws.Range("D22").Value = "'8.83"
$ws.Range("E22").Value = "  -2.38%  "

$ws.Range("D23").Value = "'5.29"
$ws.Range("E23").Value = "  +1.98%  "

$ws.Range("D24").Value = "'7.27"
$ws.Range("E24").Value = "  -2.22%  "

$ws.Range("D25").Value = "'12.51"
$ws.Range("E25").Value = "  +4.87%  "

$ws.Range("D26").Value = "'5.07"
$ws.Range("E26").Value = "  -2.16%  "

$ws.Range("D27").Value = "3.440.25"
$ws.Range("E27").Value = "  +1.35%  "

$ws.Range("D28").Value = "'75.93"
$ws.Range("E28").Value = "  -3.15%  "

$ws.Range("D29").Value = "'0.0000128"
$ws.Range("E29").Value = "  +3.16%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("E31").Value = "  +10.43%  "

$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").Value = "'8.74"
$ws.Range("E33").Value = "  -4.38%  "

$ws.Range("D34").Value = "'539.34"
$ws.Range("E34").Value = "  -4.92%  "

$ws.Range("D35").Value = "'1.41"
$ws.Range("E35").Value = "  -5.38%  "

$ws.Range("D36").Value = "'1.93"
$ws.Range("E36").Value = "  -2.89%  "

$ws.Range("D37").Value = "'6.94"
$ws.Range("E37").Value = "  +12.24%  "

$ws.Range("D38").Value = "'0.136"
$ws.Range("E38").Value = "  -11.87%  "

$ws.Range("D39").Value = "'22.32"
$ws.Range("E39").Value = "  -3.48%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").Value = "'21.57"
$ws.Range("E41").Value = "  +3.12%  "

$ws.Range("D42").Value = "'0.390"
$ws.Range("E42").Value = "  -4.13%  "

$ws.Range("D43").Value = "'1.98"
$ws.Range("E43").Value = "  -2.92%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'155.98"
$ws.Range("E45").Value = "  -2.64%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.90"
$ws.Range("E46").Value = "  -4.48%  "

$ws.Range("D47").Value = "'178.25"
$ws.Range("E47").Value = "  -5.24%  "

$ws.Range("D48").Value = "'44.48"
$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("D49").Value = "'1.29"
$ws.Range("E49").Value = "  -2.39%  "

$ws.Range("D50").Value = "'4.20"
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("D51").Value = "'0.624"
$ws.Range("E51").Value = "  -1.78%  "
